$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the header text (capitalisation change)
$ws.Range("A1").Value = "Cluster name"

# Rename an existing cluster (aged care home renamed)
$ws.Range("A6").Value = "3975 Aurrum Aged Care Brunswick West"

# Remove rows that no longer appear in the published list (delete bottom-up
# so row numbers of earlier rows stay valid while deleting)
$ws.Rows.Item(33).Delete()  # Warburton Lodge Warburton
$ws.Rows.Item(32).Delete()  # Sunny Ridge Strawberry Farm Main Ridge
$ws.Rows.Item(28).Delete()  # St Brigid's Parish Primary School Mordialloc
$ws.Rows.Item(26).Delete()  # JBS Australia Brooklyn
$ws.Rows.Item(17).Delete()  # 44893 Greenhills Primary School Greensborough
$ws.Rows.Item(14).Delete()  # 44631 Mount Evelyn Primary School
$ws.Rows.Item(12).Delete()  # 44593 Torquay P-6 College Torquay

# Insert the newly added cluster in its correct alphabetical position
# (between "Confirmed Omicron Variant..." and "Green Gables Lodge Warburton")
$ws.Rows.Item(22).Insert()
$ws.Range("A22").Value = "Feathertop Chalet Harrietville"
$ws.Range("B22").Value = 14

# Update the "Active cases" counts for every remaining cluster row
$ws.Range("B2").Value = 28
$ws.Range("B3").Value = 13
$ws.Range("B4").Value = 20
$ws.Range("B5").Value = 12
$ws.Range("B6").Value = 11
$ws.Range("B7").Value = 16
$ws.Range("B8").Value = 14
$ws.Range("B9").Value = 12
$ws.Range("B10").Value = 12
$ws.Range("B11").Value = 20
$ws.Range("B12").Value = 13
$ws.Range("B13").Value = 17
$ws.Range("B14").Value = 16
$ws.Range("B15").Value = 11
$ws.Range("B16").Value = 12
$ws.Range("B17").Value = 12
$ws.Range("B18").Value = 36
$ws.Range("B19").Value = 32
$ws.Range("B20").Value = 14
$ws.Range("B21").Value = 14
$ws.Range("B23").Value = 24
$ws.Range("B24").Value = 10
$ws.Range("B25").Value = 45
$ws.Range("B26").Value = 18
$ws.Range("B27").Value = 17
